$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1807.2222
$ws.Cells.Item(28, 9).Value = 1807.2222
$ws.Cells.Item(28, 11).Value = 1807.2222
$ws.Cells.Item(28, 13).Value = -1322.2222
$ws.Cells.Item(106, 8).Value = 10060.692
$ws.Cells.Item(106, 9).Value = 10649.083
$ws.Cells.Item(106, 11).Value = 10649.083
$ws.Cells.Item(106, 13).Value = -10018.083
$ws.Cells.Item(132, 8).Value = 6334.9697
$ws.Cells.Item(132, 9).Value = 6848.143
$ws.Cells.Item(132, 10).Value = 3461.2
$ws.Cells.Item(132, 11).Value = 20544.429
$ws.Cells.Item(132, 12).Value = 10383.6
$ws.Cells.Item(132, 13).Value = -18014.429
$ws.Cells.Item(132, 14).Value = -15443.6
$ws.Cells.Item(138, 8).Value = 519350.6
$ws.Cells.Item(138, 9).Value = 2170.5715
$ws.Cells.Item(138, 10).Value = 566366.9399999999
$ws.Cells.Item(138, 11).Value = 6511.7145
$ws.Cells.Item(138, 12).Value = 1699100.82
$ws.Cells.Item(138, 13).Value = -1371.7145
$ws.Cells.Item(138, 14).Value = -1709380.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3746.6667
$ws.Cells.Item(32, 9).Value = 3551.2969
$ws.Cells.Item(32, 10).Value = 9998.5
$ws.Cells.Item(32, 11).Value = 3551.2969
$ws.Cells.Item(32, 12).Value = 9998.5
$ws.Cells.Item(32, 13).Value = -3264.2969
$ws.Cells.Item(32, 14).Value = -10572.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 203
$ws.Cells.Item(7, 9).Value = 203
$ws.Cells.Item(7, 11).Value = 203
$ws.Cells.Item(7, 13).Value = -90
$ws.Cells.Item(22, 8).Value = 470.2
$ws.Cells.Item(22, 9).Value = 462.75
$ws.Cells.Item(22, 11).Value = 462.75
$ws.Cells.Item(22, 13).Value = -289.75
$ws.Cells.Item(86, 8).Value = 2047.4642
$ws.Cells.Item(86, 9).Value = 2148
$ws.Cells.Item(86, 10).Value = 1946.9286
$ws.Cells.Item(86, 11).Value = 2148
$ws.Cells.Item(86, 12).Value = 1946.9286
$ws.Cells.Item(86, 13).Value = -1025
$ws.Cells.Item(86, 14).Value = -4192.9286
$ws.Cells.Item(89, 8).Value = 2047.4642
$ws.Cells.Item(89, 9).Value = 2148
$ws.Cells.Item(89, 10).Value = 1946.9286
$ws.Cells.Item(89, 11).Value = 10740
$ws.Cells.Item(89, 12).Value = 9734.643
$ws.Cells.Item(89, 13).Value = -5124
$ws.Cells.Item(89, 14).Value = -20966.643

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 314.9
$ws.Cells.Item(7, 9).Value = 129.4
$ws.Cells.Item(7, 10).Value = 500.4
$ws.Cells.Item(7, 11).Value = 129.4
$ws.Cells.Item(7, 12).Value = 500.4
$ws.Cells.Item(7, 13).Value = -16.40000000000001
$ws.Cells.Item(7, 14).Value = -726.4
$ws.Cells.Item(12, 8).Value = 2857383.5
$ws.Cells.Item(12, 9).Value = 152.5
$ws.Cells.Item(12, 10).Value = 4000276
$ws.Cells.Item(12, 11).Value = 152.5
$ws.Cells.Item(12, 12).Value = 4000276
$ws.Cells.Item(12, 13).Value = 17.5
$ws.Cells.Item(12, 14).Value = -4000616
$ws.Cells.Item(22, 8).Value = 100196
$ws.Cells.Item(22, 9).Value = 185.25
$ws.Cells.Item(22, 10).Value = 233543.67
$ws.Cells.Item(22, 11).Value = 185.25
$ws.Cells.Item(22, 12).Value = 233543.67
$ws.Cells.Item(22, 13).Value = 164.75
$ws.Cells.Item(22, 14).Value = -234243.67
$ws.Cells.Item(31, 8).Value = 1370.4038
$ws.Cells.Item(31, 9).Value = 1379.0435
$ws.Cells.Item(31, 10).Value = 1363.5518
$ws.Cells.Item(31, 11).Value = 1379.0435
$ws.Cells.Item(31, 12).Value = 1363.5518
$ws.Cells.Item(31, 13).Value = -1084.0435
$ws.Cells.Item(31, 14).Value = -1953.5518
$ws.Cells.Item(34, 8).Value = 1370.4038
$ws.Cells.Item(34, 9).Value = 1379.0435
$ws.Cells.Item(34, 10).Value = 1363.5518
$ws.Cells.Item(34, 11).Value = 1379.0435
$ws.Cells.Item(34, 12).Value = 1363.5518
$ws.Cells.Item(34, 13).Value = -1177.0435
$ws.Cells.Item(34, 14).Value = -1767.5518

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 629.6429000000001
$ws.Cells.Item(113, 9).Value = 579.2857
$ws.Cells.Item(113, 11).Value = 1737.8571
$ws.Cells.Item(113, 13).Value = 432.1428999999998
$ws.Cells.Item(131, 8).Value = 52634120
$ws.Cells.Item(131, 10).Value = 3135.2666
$ws.Cells.Item(131, 12).Value = 9405.799800000001
$ws.Cells.Item(131, 14).Value = -19485.7998
$ws.Cells.Item(139, 8).Value = 1610.0769
$ws.Cells.Item(139, 9).Value = 1560.72
$ws.Cells.Item(139, 10).Value = 1698.2142
$ws.Cells.Item(139, 11).Value = 4682.16
$ws.Cells.Item(139, 12).Value = 5094.642599999999
$ws.Cells.Item(139, 13).Value = 457.8400000000001
$ws.Cells.Item(139, 14).Value = -15374.6426

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 300
$ws.Cells.Item(2, 9).Value = 317.16666
$ws.Cells.Item(2, 10).Value = 282.83334
$ws.Cells.Item(2, 11).Value = 317.16666
$ws.Cells.Item(2, 12).Value = 282.83334
$ws.Cells.Item(2, 13).Value = -204.16666
$ws.Cells.Item(2, 14).Value = -508.83334
$ws.Cells.Item(14, 8).Value = 1500
$ws.Cells.Item(14, 9).Value = 2500
$ws.Cells.Item(14, 10).Value = 500
$ws.Cells.Item(14, 11).Value = 2500
$ws.Cells.Item(14, 12).Value = 500
$ws.Cells.Item(14, 13).Value = -2332
$ws.Cells.Item(14, 14).Value = -836
$ws.Cells.Item(122, 8).Value = 2560
$ws.Cells.Item(122, 9).Value = 2834.5557
$ws.Cells.Item(122, 11).Value = 8503.667099999999
$ws.Cells.Item(122, 13).Value = -6053.667099999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1298.3334
$ws.Cells.Item(7, 9).Value = 1340
$ws.Cells.Item(7, 10).Value = 1246.25
$ws.Cells.Item(7, 11).Value = 1340
$ws.Cells.Item(7, 12).Value = 1246.25
$ws.Cells.Item(7, 13).Value = -1228
$ws.Cells.Item(7, 14).Value = -1470.25
$ws.Cells.Item(22, 8).Value = 870
$ws.Cells.Item(22, 9).Value = 375
$ws.Cells.Item(22, 10).Value = 1200
$ws.Cells.Item(22, 11).Value = 375
$ws.Cells.Item(22, 12).Value = 1200
$ws.Cells.Item(22, 13).Value = -80
$ws.Cells.Item(22, 14).Value = -1790
$ws.Cells.Item(25, 8).Value = 2605.6667
$ws.Cells.Item(25, 9).Value = 1200
$ws.Cells.Item(25, 10).Value = 3007.2856
$ws.Cells.Item(25, 11).Value = 1200
$ws.Cells.Item(25, 12).Value = 3007.2856
$ws.Cells.Item(25, 13).Value = -970
$ws.Cells.Item(25, 14).Value = -3467.2856
$ws.Cells.Item(27, 8).Value = 870
$ws.Cells.Item(27, 9).Value = 375
$ws.Cells.Item(27, 10).Value = 1200
$ws.Cells.Item(27, 11).Value = 375
$ws.Cells.Item(27, 12).Value = 1200
$ws.Cells.Item(27, 13).Value = -268
$ws.Cells.Item(27, 14).Value = -1414
$ws.Cells.Item(119, 8).Value = 24000
$ws.Cells.Item(119, 9).Value = 20000
$ws.Cells.Item(119, 11).Value = 20000
$ws.Cells.Item(119, 13).Value = -15162
$ws.Cells.Item(126, 8).Value = 1298.3334
$ws.Cells.Item(126, 9).Value = 1340
$ws.Cells.Item(126, 10).Value = 1246.25
$ws.Cells.Item(126, 11).Value = 4020
$ws.Cells.Item(126, 12).Value = 3738.75
$ws.Cells.Item(126, 13).Value = -1550
$ws.Cells.Item(126, 14).Value = -8678.75
$ws.Cells.Item(132, 8).Value = 2510.3794
$ws.Cells.Item(132, 9).Value = 2144.611
$ws.Cells.Item(132, 11).Value = 6433.833
$ws.Cells.Item(132, 13).Value = -3903.833
$ws.Cells.Item(136, 8).Value = 2030.25
$ws.Cells.Item(136, 9).Value = 1448.4
$ws.Cells.Item(136, 11).Value = 4345.200000000001
$ws.Cells.Item(136, 13).Value = -1795.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1822.7297
$ws.Cells.Item(132, 9).Value = 1582.5172
$ws.Cells.Item(132, 10).Value = 2693.5
$ws.Cells.Item(132, 11).Value = 4747.5516
$ws.Cells.Item(132, 12).Value = 8080.5
$ws.Cells.Item(132, 13).Value = -2217.5516
$ws.Cells.Item(132, 14).Value = -13140.5
$ws.Cells.Item(136, 8).Value = 1244.1852
$ws.Cells.Item(136, 9).Value = 1195.1904
$ws.Cells.Item(136, 11).Value = 3585.5712
$ws.Cells.Item(136, 13).Value = -1035.5712
